$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of user-story data to append below the existing header/first story
# (Row 1 = header "User Stories"/"Points", Row 2 = first story already present)
$rows = @(
    @{ Row = 3; B = "As the system, I reject registration attempts for usernames that already exist"; C = $null },
    @{ Row = 4; B = "1 point"; C = 1 },
    @{ Row = 5; B = "As a user, I can login with a username and password"; C = $null },
    @{ Row = 6; B = "2 points"; C = 2 },
    @{ Row = 7; B = "As the system, I reject login attempts with invalid credentials"; C = $null },
    @{ Row = 8; B = "1 point"; C = 1 },
    @{ Row = 9; B = "As a customer, I can apply for a new bank account with a starting balance."; C = $null },
    @{ Row = 10; B = "3 points"; C = 3 },
    @{ Row = 11; B = "As a customer, I can view the balance of a specific account"; C = $null },
    @{ Row = 12; B = "1 point"; C = 1 },
    @{ Row = 13; B = "As a customer, I can make a deposit to a specific account"; C = $null },
    @{ Row = 14; B = "2 points"; C = 2 },
    @{ Row = 15; B = "As a customer, I can make a withdrawal from a specific account"; C = $null },
    @{ Row = 16; B = "2 points"; C = 2 },
    @{ Row = 17; B = "As the system, I reject and prevent overdrafts"; C = $null },
    @{ Row = 18; B = "1 points"; C = 1 },
    @{ Row = 19; B = "As the system, I reject deposits or withdrawals of negative money"; C = $null },
    @{ Row = 20; B = "2 points"; C = 2 },
    @{ Row = 21; B = "As the system, I reject any transactions of unapproved accounts"; C = $null },
    @{ Row = 22; B = "1 point"; C = 1 },
    @{ Row = 23; B = "As an employee, I can approve or reject an account."; C = $null },
    @{ Row = 24; B = "2 points"; C = 2 },
    @{ Row = 25; B = "As a customer, I can post a money transfer to another account."; C = $null },
    @{ Row = 26; B = "3 points"; C = 3 },
    @{ Row = 27; B = "As the system, I reject invalid transfers (negative amounts or overdrafts)"; C = $null },
    @{ Row = 28; B = "1 point"; C = 1 },
    @{ Row = 29; B = "As an employee, I can view a log of all transactions."; C = $null },
    @{ Row = 30; B = "2 points"; C = 2 },
    @{ Row = 31; B = "File I/O integration tests"; C = $null },
    @{ Row = 32; B = "10 points"; C = 10 },
    @{ Row = 33; B = "Database integration tests"; C = $null },
    @{ Row = 34; B = "NOTE: database tests are ignored for scoring purposes, but you should still run the tests locally to help you validate your JDBC code"; C = $null }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    if ($null -ne $r.C) {
        $ws.Cells.Item($r.Row, 3).Value = $r.C
    }
}

# Widen column B to fit the longer text
$ws.Columns.Item(2).ColumnWidth = 74.75

# Update the page setup to match the saved workbook (portrait orientation)
$ws.PageSetup.Orientation = 1

# Selection state as captured when the file was last saved
$ws.Range("C32").Select()
